$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.339.30"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "3.362.90"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'182.42"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "'537.58"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "'0.601"
$ws.Range("E7").Value = "  -1.33%  "
$ws.Range("D8").Value = "3.354.81"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'0.625"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").Value = "'55.04"
$ws.Range("E11").Value = "  -7.97%  "
$ws.Range("E12").Value = "  +3.24%  "
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "'9.19"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").Value = "3.909.65"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.120"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.366.02"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "'17.97"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").Value = "65.526.73"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "'0.981"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "'389.90"
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("D23").Value = "'11.87"
$ws.Range("E23").Value = "  +3.10%  "
$ws.Range("D24").Value = "'4.21"
$ws.Range("E24").Value = "  +6.00%  "
$ws.Range("D25").Value = "'82.85"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Value = "'3.78"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").Value = "'6.11"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +3.79%  "
$ws.Range("D29").Value = "'11.55"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").Value = "'8.44"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "'29.45"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "'663.41"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value = "'6.78"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("D34").Value = "'11.42"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "'57.77"
$ws.Range("E36").Value = "  -3.58%  "
$ws.Range("D37").Value = "'37.66"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'0.397"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").Value = "0.0₃0774"
$ws.Range("E40").Value = "  +8.14%  "
$ws.Range("D41").Value = "'2.77"
$ws.Range("E41").Value = "  +7.16%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'3.28"
$ws.Range("E42").Value = "  +15.13%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "'0.128"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").Value = "3.009.64"
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("D46").Value = "'2.75"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").Value = "'0.0410"
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("D49").Value = "'3.17"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'139.29"
$ws.Range("E51").Value = "  +2.36%  "
